# Journal de bord CPNVoiturage - actualisation de la doc
# Fills in the journal rows 40-49 with new entries, and extends the
# bottom of the table by two rows (table now spans C2:F90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the new journal entries (rows 40-49) -----------------
# Strings are written in a specific order so that the underlying shared
# string table is built in the same sequence as the source workbook.

$ws.Range("C41").Value = "Remplissage de la base de données"
$ws.Range("F41").Value = "Ajout de données de test"
$ws.Range("C40").Value = "Sprint 1 review + site en ligne"
$ws.Range("F40").Value = "Accès au site en ligne (mycpnv) et cloture du sprint 2"
$ws.Range("C42").Value = "Création du DB connetor"
$ws.Range("F42").Value = "Tests avec la connexion à la DB et des fonctions de requêtes select et insert"
$ws.Range("C43").Value = "Affichage des villes sur les register"
$ws.Range("C44").Value = "Gestion du formulaire de register"
$ws.Range("F44").Value = "Formulaire HTML et vérification des données coté PHP"
$ws.Range("F45").Value = "Requête et hash du password"
$ws.Range("C46").Value = "Gestion du formulaire de login"
$ws.Range("F46").Value = "Formulaire HTML et vérification des données coté PHP"
$ws.Range("C45").Value = "Requêtes du register"
$ws.Range("C47").Value = "Requêtes du login"
$ws.Range("F47").Value = "Requête et password verify"
$ws.Range("C49").Value = "Création de la requête pour l'affichage des voitures"
$ws.Range("F49").Value = "Création d'une requête unique pour afficher toutes les voitures du jour"
$ws.Range("C48").Value = "Sécurisation et gestion des erreurs du login"
$ws.Range("F48").Value = "Affichage des erreurs et réaffichage des données envoyées en cas d'erreur"

# Dates (column D) and durations in minutes (column E)
$ws.Range("D40").Value = 44263
$ws.Range("E40").Value = 60

$ws.Range("D41").Value = 44263
$ws.Range("E41").Value = 80

$ws.Range("D42").Value = 44266
$ws.Range("E42").Value = 60

$ws.Range("D43").Value = 44266
$ws.Range("E43").Value = 15

$ws.Range("D44").Value = 44266
$ws.Range("E44").Value = 30

$ws.Range("D45").Value = 44266
$ws.Range("E45").Value = 30

$ws.Range("D46").Value = 44266
$ws.Range("E46").Value = 20

$ws.Range("D47").Value = 44266
$ws.Range("E47").Value = 20

$ws.Range("D48").Value = 44267
$ws.Range("E48").Value = 30

$ws.Range("D49").Value = 44267
$ws.Range("E49").Value = 60

# --- 2. Extend the table with two additional blank rows ---------------
# Row 88 becomes a normal blank data row (same style as the rows above).
$ws.Range("C87:F87").Copy()
$ws.Range("C88:F88").PasteSpecial(-4122)

# Row 89 used to be the "last row" of the table (heavier bottom border);
# it now becomes a normal blank data row, and that special bottom-border
# style moves down to the new last row (row 91).
$ws.Range("C89:F89").Copy()
$ws.Range("C91:F91").PasteSpecial(-4122)

$ws.Range("C87:F87").Copy()
$ws.Range("C89:F89").PasteSpecial(-4122)

# Row 90 is a near-empty spacer row, only E90 carries a style (the
# "minutes" number format, centered vertically, no border).
$ws.Range("E90").NumberFormat = '0\ "minutes"'
$ws.Range("E90").VerticalAlignment = -4108

$ws.Application.CutCopyMode = $false

# --- 3. Resize the table / autofilter to cover the new rows -----------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("C2:F90"))

# --- 4. Update the sheet view to reflect where the author left off ----
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F58").Select()
